$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016427182986297
$ws.Range("D2").Value = 1.048389751436182
$ws.Range("E2").Value = 1.017954707410919
$ws.Range("F2").Value = 1.04952443735286
$ws.Range("I2").Value = 1.038666745281588
$ws.Range("J2").Value = 1.021647055083853
$ws.Range("K2").Value = 1.051149731471575
$ws.Range("L2").Value = 1.020802360218753
$ws.Range("M2").Value = 1.052281255409021
$ws.Range("N2").Value = 1.023097911086627

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017575060847816
$ws.Range("D3").Value = 1.049190011978177
$ws.Range("E3").Value = 1.01893224704502
$ws.Range("F3").Value = 1.050595856026576
$ws.Range("I3").Value = 1.038871032416987
$ws.Range("J3").Value = 1.022429335216581
$ws.Range("K3").Value = 1.051762058123186
$ws.Range("L3").Value = 1.021585132180153
$ws.Range("M3").Value = 1.05316426841827
$ws.Range("N3").Value = 1.023881302146872

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018316930754252
$ws.Range("D4").Value = 1.049703551887514
$ws.Range("E4").Value = 1.019564426710736
$ws.Range("F4").Value = 1.051284504423114
$ws.Range("I4").Value = 1.038999498912891
$ws.Range("J4").Value = 1.022934161810573
$ws.Range("K4").Value = 1.052153348614996
$ws.Range("L4").Value = 1.022090678415468
$ws.Range("M4").Value = 1.053730420915254
$ws.Range("N4").Value = 1.024386845652532

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018628603075056
$ws.Range("D5").Value = 1.049918418405833
$ws.Range("E5").Value = 1.019830111084538
$ws.Range("F5").Value = 1.051572903920116
$ws.Range("I5").Value = 1.039052614827649
$ws.Range("J5").Value = 1.02314606590593
$ws.Range("K5").Value = 1.052316667261744
$ws.Range("L5").Value = 1.02230298105638
$ws.Range("M5").Value = 1.053967182649937
$ws.Range("N5").Value = 1.02459905067601

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018680921978472
$ws.Range("D6").Value = 1.049954435277478
$ws.Range("E6").Value = 1.019874715757676
$ws.Range("F6").Value = 1.051621262483013
$ws.Range("I6").Value = 1.039061480932479
$ws.Range("J6").Value = 1.023181626500147
$ws.Range("K6").Value = 1.052344019971794
$ws.Range("L6").Value = 1.022338614195737
$ws.Range("M6").Value = 1.054006862753339
$ws.Range("N6").Value = 1.024634661770351

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018321096157254
$ws.Range("D7").Value = 1.049706426975588
$ws.Range("E7").Value = 1.019567977127408
$ws.Range("F7").Value = 1.051288362386062
$ws.Range("I7").Value = 1.039000212153894
$ws.Range("J7").Value = 1.022936994559628
$ws.Range("K7").Value = 1.052155535524018
$ws.Range("L7").Value = 1.022093516112753
$ws.Range("M7").Value = 1.053733589444506
$ws.Range("N7").Value = 1.024389682424417

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016815297786094
$ws.Range("D8").Value = 1.048661089625061
$ws.Range("E8").Value = 1.018285145239922
$ws.Range("F8").Value = 1.049887487511467
$ws.Range("I8").Value = 1.038736555182211
$ws.Range("J8").Value = 1.021911712617704
$ws.Range("K8").Value = 1.051357689709241
$ws.Range("L8").Value = 1.021067100802874
$ws.Range("M8").Value = 1.052580753799402
$ws.Range("N8").Value = 1.023362944464533

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.01415503598944
$ws.Range("D9").Value = 1.046786289862826
$ws.Range("E9").Value = 1.016021885039992
$ws.Range("F9").Value = 1.047383470577058
$ws.Range("I9").Value = 1.038243482413679
$ws.Range("J9").Value = 1.020094572119232
$ws.Range("K9").Value = 1.049914078242586
$ws.Range("L9").Value = 1.019251041589998
$ws.Range("M9").Value = 1.050509355757034
$ws.Range("N9").Value = 1.021543223418124

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012376797610093
$ws.Range("D10").Value = 1.045514388541736
$ws.Range("E10").Value = 1.014511141304511
$ws.Range("F10").Value = 1.045690205154306
$ws.Range("I10").Value = 1.037895645213857
$ws.Range("J10").Value = 1.018876049379394
$ws.Range("K10").Value = 1.048926334229707
$ws.Range("L10").Value = 1.018035315920704
$ws.Range("M10").Value = 1.049101535318999
$ws.Range("N10").Value = 1.02032297023623

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011605649508877
$ws.Range("D11").Value = 1.044958411924014
$ws.Range("E11").Value = 1.013856506635777
$ws.Range("F11").Value = 1.044951311789371
$ws.Range("I11").Value = 1.037740495179874
$ws.Range("J11").Value = 1.018346717147492
$ws.Range("K11").Value = 1.048492619796586
$ws.Range("L11").Value = 1.017507688898918
$ws.Range("M11").Value = 1.048485545487177
$ws.Range("N11").Value = 1.019792886291843

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011319033927121
$ws.Range("D12").Value = 1.044751110323323
$ws.Range("E12").Value = 1.013613273755357
$ws.Range("F12").Value = 1.044675995894035
$ws.Range("I12").Value = 1.037682184344046
$ws.Range("J12").Value = 1.018149842106537
$ws.Range("K12").Value = 1.048330614834804
$ws.Range("L12").Value = 1.01731152162977
$ws.Range("M12").Value = 1.048255777427251
$ws.Range("N12").Value = 1.019595731665748

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011380521988614
$ws.Range("D13").Value = 1.04479561285121
$ws.Range("E13").Value = 1.013665451333946
$ws.Range("F13").Value = 1.044735090946591
$ws.Range("I13").Value = 1.037694723042759
$ws.Range("J13").Value = 1.01819208414731
$ws.Range("K13").Value = 1.048365406370622
$ws.Range("L13").Value = 1.017353608482663
$ws.Range("M13").Value = 1.048305107012088
$ws.Range("N13").Value = 1.019638033695063

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01158196139877
$ws.Range("D14").Value = 1.044941292370403
$ws.Range("E14").Value = 1.013836402415584
$ws.Range("F14").Value = 1.044928571606085
$ws.Range("I14").Value = 1.037735689084559
$ws.Range("J14").Value = 1.018330448656593
$ws.Range("K14").Value = 1.048479246862608
$ws.Range("L14").Value = 1.017491477394397
$ws.Range("M14").Value = 1.048466572420673
$ws.Range("N14").Value = 1.01977659469782

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011706051352975
$ws.Range("D15").Value = 1.045030946006555
$ws.Range("E15").Value = 1.013941721380301
$ws.Range("F15").Value = 1.045047667673156
$ws.Range("I15").Value = 1.037760839347819
$ws.Range("J15").Value = 1.018415665432846
$ws.Range("K15").Value = 1.048549267927591
$ws.Range("L15").Value = 1.01757639867754
$ws.Range("M15").Value = 1.048565929069641
$ws.Range("N15").Value = 1.01986193249167

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.012427951519364
$ws.Range("D16").Value = 1.045551176532851
$ws.Range("E16").Value = 1.014554577174231
$ws.Range("F16").Value = 1.045739122858873
$ws.Range("I16").Value = 1.03790584648219
$ws.Range("J16").Value = 1.018911143391592
$ws.Range("K16").Value = 1.048954991667732
$ws.Range("L16").Value = 1.018070307185995
$ws.Range("M16").Value = 1.049142281704191
$ws.Range("N16").Value = 1.02035811408595

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012880468200438
$ws.Range("D17").Value = 1.045876101441905
$ws.Range("E17").Value = 1.014938877900143
$ws.Range("F17").Value = 1.04617132754816
$ws.Range("I17").Value = 1.03799559171409
$ws.Range("J17").Value = 1.019221486415713
$ws.Range("K17").Value = 1.049207880843299
$ws.Range("L17").Value = 1.018379798116014
$ws.Range("M17").Value = 1.049502099653171
$ws.Range("N17").Value = 1.020668897832761

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01314430186988
$ws.Range("D18").Value = 1.046065119599636
$ws.Range("E18").Value = 1.015162988505591
$ws.Range("F18").Value = 1.046422875570972
$ws.Range("I18").Value = 1.038047501106282
$ws.Range("J18").Value = 1.019402339960163
$ws.Range("K18").Value = 1.049354806524378
$ws.Range("L18").Value = 1.018560202307227
$ws.Range("M18").Value = 1.049711358520283
$ws.Range("N18").Value = 1.02085000820999

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.013234243474581
$ws.Range("D19").Value = 1.046129484258916
$ws.Range("E19").Value = 1.01523939675313
$ws.Range("F19").Value = 1.046508553768431
$ws.Range("I19").Value = 1.038065126660939
$ws.Range("J19").Value = 1.019463978522597
$ws.Range("K19").Value = 1.04940480596757
$ws.Range("L19").Value = 1.018621695731046
$ws.Range("M19").Value = 1.049782605794168
$ws.Range("N19").Value = 1.020911734306252

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012831929012738
$ws.Range("D20").Value = 1.045841292318954
$ws.Range("E20").Value = 1.014897650787199
$ws.Range("F20").Value = 1.046125012944362
$ws.Range("I20").Value = 1.037986008151355
$ws.Range("J20").Value = 1.019188206551323
$ws.Range("K20").Value = 1.049180808255103
$ws.Range("L20").Value = 1.018346604742526
$ws.Range("M20").Value = 1.049463558407618
$ws.Range("N20").Value = 1.020635570707148

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011522647416815
$ws.Range("D21").Value = 1.044898415117351
$ws.Range("E21").Value = 1.01378606359346
$ws.Range("F21").Value = 1.044871620075269
$ws.Range("I21").Value = 1.03772364441114
$ws.Range("J21").Value = 1.018289710877292
$ws.Range("K21").Value = 1.048445748637298
$ws.Range("L21").Value = 1.017450883501681
$ws.Range("M21").Value = 1.048419051449932
$ws.Range("N21").Value = 1.019735799066201

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010698425588204
$ws.Range("D22").Value = 1.044301036299226
$ws.Range("E22").Value = 1.013086744894684
$ws.Range("F22").Value = 1.044078597054144
$ws.Range("I22").Value = 1.037554744627711
$ws.Range("J22").Value = 1.017723300497727
$ws.Range("K22").Value = 1.047978356260124
$ws.Range("L22").Value = 1.016886647943687
$ws.Range("M22").Value = 1.047756761807001
$ws.Range("N22").Value = 1.019168584318928

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.011135459201395
$ws.Range("D23").Value = 1.04461815007401
$ws.Range("E23").Value = 1.013457507234019
$ws.Range("F23").Value = 1.044499464879126
$ws.Range("I23").Value = 1.037644655201384
$ws.Range("J23").Value = 1.018023707152083
$ws.Range("K23").Value = 1.048226625849023
$ws.Range("L23").Value = 1.017185860831585
$ws.Range("M23").Value = 1.048108382220434
$ws.Range("N23").Value = 1.019469417585191

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012853862127242
$ws.Range("D24").Value = 1.045857022623829
$ws.Range("E24").Value = 1.014916279685586
$ws.Range("F24").Value = 1.046145942218525
$ws.Range("I24").Value = 1.037990339902626
$ws.Range("J24").Value = 1.019203244797149
$ws.Range("K24").Value = 1.049193042985439
$ws.Range("L24").Value = 1.018361603759597
$ws.Range("M24").Value = 1.049480975443441
$ws.Range("N24").Value = 1.020650630309007

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014843601981001
$ws.Range("D25").Value = 1.047274854250745
$ws.Range("E25").Value = 1.016607323215551
$ws.Range("F25").Value = 1.048035030547581
$ws.Range("I25").Value = 1.03837432536248
$ws.Range("J25").Value = 1.02056559248068
$ws.Range("K25").Value = 1.050291752683983
$ws.Range("L25").Value = 1.019721416963889
$ws.Range("M25").Value = 1.051049598612692
$ws.Range("N25").Value = 1.022014912682513
